# Insert a new data row at row 308 ("Feria Lagunitas de Puerto Montt" -
# Zanahoria subset), pushing the existing rows 308-360 down to 309-361.
# This mirrors the commit's weekly price-update: a brand-new observation
# was added into the middle of the date-sorted block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 308:360 down to 309:361, leaving row 308 blank for the new record.
$ws.Rows.Item(308).Insert()

# Populate the new row 308 with the new observation.
$ws.Range("A308").Value = 4
$ws.Range("B308").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C308").Value = "Los Lagos"
$ws.Range("D308").Value = 44694
$ws.Range("D308").NumberFormat = $ws.Range("D309").NumberFormat
$ws.Range("E308").Value = 10
$ws.Range("F308").Value = 100114013
$ws.Range("G308").Value = "Zanahoria"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 900
$ws.Range("K308").Value = 7500
$ws.Range("L308").Value = 8000
$ws.Range("M308").Value = 7750
$ws.Range("N308").Value = "`$/saco 20 kilos"
$ws.Range("O308").Value = "Provincia de Llanquihue"
$ws.Range("P308").Value = 388
$ws.Range("Q308").Value = 20
$ws.Range("R308").Value = "Hortaliza"
